# edit.ps1 - applies the commit "Added a second paragraph on methods and
# ATRIP, also finished last paragraph off" to The Report .docx

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1) "After all players have had their turn..." -> add comma after "all",
#    and "so that the game can continue playing" -> "to continue playing,"
Replace-Text "After all players have had their turn, the method calls upon itself at the end so that the game can continue playing and players can have subsequent turns." "After all, players have had their turn, the method calls upon itself at the end to continue playing, and players can have subsequent turns."

# 2) "When the player enters the weight of the pebble they wish to remove" -> "the pebble's weight"
#    and "in player's ArrayList" -> "in the player's ArrayList"
Replace-Text "When the player enters the weight of the pebble they wish to remove, it calls upon a method that removes the first instance of that integer in player's ArrayList of pebbles and" "When the player enters the pebble's weight they wish to remove, it calls upon a method that removes the first instance of that integer in the player's ArrayList of pebbles and"

# 3) "was set to true at any point" -> "were set to true at any point"
$q = [char]34
$find3 = "The program was structured so that if " + $q + "finished" + $q + " was set to true at any point"
$repl3 = "The program was structured so that if " + $q + "finished" + $q + " were set to true at any point"
Replace-Text $find3 $repl3

# 4) Remove comma after "strictly positive number of players"
Replace-Text "Finally, to ensure that the game has a strictly positive number of players, when asked for the number playing" "Finally, to ensure that the game has a strictly positive number of players when asked for the number playing"

# 5) "pebbles could be removed, but errors were being caught." -> "pebbles could be removed but caught errors."
Replace-Text "pebbles could be removed, but errors were being caught." "pebbles could be removed but caught errors."

# 6) PebbleGameTest paragraph rewording
Replace-Text "The next test method was on writing to a text file when the player discards a pebble. Two tests were also carried out here using values that were valid. In a similar fashion to the previous method, the first test was carried out using the A white bag, and the second used the C white bag such that both boundary cases were tested. Again, the second test here was also designed to check that the file appended the text from the first test, instead of overwriting it." "The following test method was on writing to a text file when the player discards a pebble. Two tests were also carried out here using valid values. In a similar fashion to the previous method, the first test was carried out using the A white bag, and the second used the C white bag such that both boundary cases were tested. Again, the second test was also designed to check that the file appended the text from the first test instead of overwriting it."

# 7) "three separate games" -> "three different games", and the following sentences reworded
Replace-Text "A CSV file with 11 pebbles was used to test three separate games, varying in player counts of 1, 4 and 100. The games were checked to see if a dummy black bag would end up with at least 11, 44 and 1100 pebbles respectively, after having read the CSV files. And as the CSV files could not contain negative numbers, or contain fewer than 11 integers, CSV files were created which separately fail those criteria, and throw the correct exceptions when read within this method." "A CSV file with 11 pebbles was used to test three different games, varying in player counts of 1, 4 and 100. After reading the CSV files, the games were checked to see if a dummy black bag would end up with at least 11, 44 and 1100 pebbles, respectively. And as the CSV files could not contain negative numbers or contain fewer than 11 integers, CSV files were created which separately failed those criteria and threw the right exceptions when read within this method."

# 8) "In order to test the method which draws a pebble..." paragraph reworded
$rsq = [char]8217
$find8 = "In order to test the method which draws a pebble from a black bag in the game, another test was designed to apply this method and assert that the player" + $rsq + "s pebble count increased by one, and that the black bag it was drawn from decreased its count by one. If both events take place, it can be surmised that the draw worked correctly. Additionally, this method was also created to check that, if a black bag was empty, a non-empty white bag would empty its contents into a black bag, and the draw would take place as normal. This was tested by asserting that the white bag was empty after the draw, and the black bag ended with one fewer pebble than the white bag started with (to account for the pebble given to the player). Similar to this test, another testing method checked that the initial draw took place correctly by asserting that the player" + $rsq + "s pebble count increased by 10, and that the black bag decreased by the same amount after this draw took place."
$repl8 = "In order to test the method which draws a pebble from a black bag in the game, another test was designed to apply this method and assert that the player's pebble count increased by one and that the black bag was drawn from decreased its count by one. If both events take place, it can be surmised that the draw worked correctly. Additionally, this method was also created to check that if a black bag were empty, a non-empty white bag would empty its contents into a black bag, and the draw would take place as expected. This was tested by asserting that the white bag was empty after the draw, and the black bag ended with one fewer pebble than the white bag started with (to account for the pebble given to the player). Like this test, another testing method checked that the initial draw took place correctly by asserting that the player's pebble count increased by 10 and that the black bag decreased by the same amount after this draw."
Replace-Text $find8 $repl8

# 9) Insert the new ATRIP paragraph before the "The file PlayerTest.java..." paragraph
$atripText = "The tests were designed to follow the ATRIP properties to ensure they were good tests. They are automatic since they do not require user input or files and automatically check that each criterion has been met. This was done with the use of the methods writer, reader and test input. The writer method would create temporary files and write to them, and the reader would read from them for validation. The method testInput would take a string and then convert this into bytes. It would then convert these bytes into a byte input stream, replacing the standard input stream before setting the scanner in the pebble game to the correct input stream. The tests were through since I assured that they tested for the required properties and that each input would catch any thrown errors. The tests met the repeatable property because I made sure they used the same input each time and did not set any checks on values that could change between tests (e.g. their value was created randomly). If there were times that I had to check a value that could change, I would check a range of values. The tests were independent of each other using the two setup methods that would reset all values before a test so that the result of another test wouldn't influence the result of another. Finally, I made sure the testing was professional but included plenty of comments to make it easier to understand, diagnose errors, and maintain both the test file and production code. Also, this was a reason for splitting up the tests for each method and using a method that would set up for each test to make the code more readable."

$anchor = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "The file PlayerTest.java was used*") {
        $anchor = $i
        break
    }
}
$d.Paragraphs($anchor).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs($anchor).Range.Text = $atripText

# 10) playerThreadTest paragraph: reword main body text
$find10 = "This test started off by checking that " + $q + "E" + $q + " ends the game, by asserting that the Boolean " + $q + "finished" + $q + " is indeed triggered by that input. Similarly, the test was made to check that a total pebble weight of 100 from the initial draw also ends the game immediately."
$repl10 = "Conducting these tests would be like testing the game if there was only one player and a single thread. This test started by checking that " + $q + "E" + $q + " ends the game by asserting that the Boolean " + $q + "finished" + $q + " is triggered by that input. Similarly, the test was made to check that a total pebble weight of 100 from the initial draw ends the game immediately."
Replace-Text $find10 $repl10

# 11) Replace the trailing bold "(TODO- ...)" note with a single (still bold) space,
#     followed by a new, non-bold sentence.
$rng = $d.Content
$todoFind = " (TODO- finish off explaining playerThreadTest() method here  from test on removing pebble and adding a new one, line 399)"
$found = $rng.Find.Execute($todoFind, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = " "
    $insPoint = $d.Range($rng.End, $rng.End)
    $sq = [char]39
    $tail = "Unlike the other tests that were checking win or end conditions, this test ensured that the playerThread was functioning as expected by inputting a 2 (an arbitrary value that would always be in the players pebbles). Then several checks could be made that the 2 had been removed and a random pebble from the bags had been inserted as expected. The last 2 tests were both done similarly to check that pebbles that were not inside the player" + $sq + "s pebbles and that values that weren" + $sq + "t integers should be thrown and caught."
    $insPoint.InsertAfter($tail)
    $insPoint.Font.Bold = 0
}

Write-Host "Done basic replacements"
